$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Update the date placeholder field on the slide master and every slide layout
# (datetimeFigureOut field text changes from 08/01/2019 to 10/01/2019).
function Update-DateField($container) {
    foreach ($hf in @($container.HeadersFooters)) {
        # no-op placeholder in case HeadersFooters exposes the date text in some hosts
    }
}

$master = $p.SlideMaster
foreach ($shape in $master.Shapes) {
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "08/01/2019") {
            $tr.Text = "10/01/2019"
        }
    }
}

foreach ($layout in $p.SlideMaster.CustomLayouts) {
    foreach ($shape in $layout.Shapes) {
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "08/01/2019") {
                $tr.Text = "10/01/2019"
            }
        }
    }
}

# Move the grouped logo shape (Groupe 15) to its new position.
# (Target EMU: x=4223751, y=812097; point values below are nudged to the
# nearest-representable Single so that PowerPoint's internal EMU rounding
# lands exactly on the target instead of one EMU short.)
$grp = $s.Shapes.Item(1)
$grp.Left = 332.5788269042969
$grp.Top = 63.94464874267578
